# Apply updates described by the commit "Update countries & provincias Spain"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (row 1, col A) ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 12:04"

# --- Swap Seychelles / Montserrat rows (205 = Seychelles, 206 = Montserrat now) ---
$ws.Range("A205").Value = "Seychelles"
$ws.Range("A206").Value = "Montserrat"

# --- Row-by-row numeric updates (country statistics refresh) ---

# Row 18 - Belgica
$ws.Range("B18").Value = 52011
$ws.Range("C18").Value = 591
$ws.Range("D18").Value = 13201
$ws.Range("E18").Value = 30289
$ws.Range("F18").Value = 508
$ws.Range("G18").Value = 106
$ws.Range("H18").Value = 8521

# Row 32 - Austria
$ws.Range("B32").Value = 16409
$ws.Range("C32").Value = 28
$ws.Range("D32").Value = 11007
$ws.Range("E32").Value = 5157
$ws.Range("F32").Value = 77
$ws.Range("G32").Value = 5
$ws.Range("H32").Value = 245

# Row 34 - Polonia
$ws.Range("B34").Value = 15774
$ws.Range("C34").Value = 22
$ws.Range("D34").Value = 13836
$ws.Range("E34").Value = 1324
$ws.Range("F34").Value = 81
$ws.Range("G34").Value = 5
$ws.Range("H34").Value = 614

# Row 37 - Banglades
$ws.Range("B37").Value = 14811
$ws.Range("C37").Value = 312
$ws.Range("D37").Value = 6423
$ws.Range("E37").Value = 7490
$ws.Range("F37").Value = 232
$ws.Range("G37").Value = 10
$ws.Range("H37").Value = 898

# Row 40 - Dinamarca
$ws.Range("B40").Value = 13112
$ws.Range("C40").Value = 336
$ws.Range("D40").Value = 2494
$ws.Range("E40").Value = 9675
$ws.Range("G40").Value = 13
$ws.Range("H40").Value = 943

# Row 49 - Chequia
$ws.Range("B49").Value = 8034
$ws.Range("C49").Value = 3
$ws.Range("D49").Value = 4372
$ws.Range("E49").Value = 3392

# Row 59 - Moldavia
$ws.Range("D59").Value = 1582
$ws.Range("E59").Value = 3140

# Row 72 - Uzbekistan
$ws.Range("B72").Value = 2314
$ws.Range("C72").Value = 16
$ws.Range("D72").Value = 1721
$ws.Range("E72").Value = 583

# Row 92 - Hong Kong
$ws.Range("D92").Value = 960
$ws.Range("E92").Value = 81

# Row 101 - Albania
$ws.Range("B101").Value = 850
$ws.Range("C101").Value = 8
$ws.Range("D101").Value = 620
$ws.Range("E101").Value = 199

# Row 114 - Georgia
$ws.Range("E114").Value = 325
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 10

# Row 205 (now Seychelles)
$ws.Range("D205").Value = 8
$ws.Range("F205").Value = 0
$ws.Range("H205").Value = 0

# Row 206 (now Montserrat)
$ws.Range("D206").Value = 7
$ws.Range("F206").Value = 1
$ws.Range("H206").Value = 1
